$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates: domain changed from "i-preproducciongestion..." to "preproducciongestion..."
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
$ws.Range("E3").Value = "'12112001793"
$ws.Range("G3").Value = "'10/04/2021"

# Column width changes
$ws.Columns.Item(5).ColumnWidth = 19.42578125
$ws.Columns.Item(7).ColumnWidth = 16.28515625

# Selection change
$ws.Range("B4").Select()
